# Update gh-pages to output generated at 456a3b4
# Applies refreshed "want to go" counts (column F) and "lowest price" (column G)
# for the 广州-漫展信息 workbook across its four sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 638
$ws.Range("G2").Value = "已售罄"
$ws.Range("F3").Value = 290
$ws.Range("F5").Value = 768
$ws.Range("G5").Value = 70
$ws.Range("F6").Value = 461
$ws.Range("F8").Value = 203
$ws.Range("F10").Value = 278
$ws.Range("F11").Value = 7161
$ws.Range("F12").Value = 77
$ws.Range("F14").Value = 1124
$ws.Range("F16").Value = 567
$ws.Range("F17").Value = 395
$ws.Range("F19").Value = 128
$ws.Range("F22").Value = 10
$ws.Range("F23").Value = 44
$ws.Range("F24").Value = 118
$ws.Range("F26").Value = 201
$ws.Range("F28").Value = 346
$ws.Range("F30").Value = 1058
$ws.Range("F32").Value = 77
$ws.Range("F33").Value = 2053
$ws.Range("F34").Value = 575
$ws.Range("F35").Value = 5
$ws.Range("F36").Value = 11
$ws.Range("F37").Value = 35

# ---- Sheet: 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 291

# ---- Sheet: 本地生活 (Local Life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 360

# ---- Sheet: 全部类型 (All Types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 360
$ws.Range("F3").Value = 638
$ws.Range("G3").Value = "已售罄"
$ws.Range("F4").Value = 290
$ws.Range("F6").Value = 768
$ws.Range("G6").Value = 70
$ws.Range("F8").Value = 461
$ws.Range("F10").Value = 203
$ws.Range("F12").Value = 278
$ws.Range("F13").Value = 7161
$ws.Range("F14").Value = 77
$ws.Range("F17").Value = 1124
$ws.Range("F19").Value = 567
$ws.Range("F20").Value = 395
$ws.Range("F23").Value = 128
$ws.Range("F25").Value = 291
$ws.Range("F29").Value = 10
$ws.Range("F30").Value = 44
$ws.Range("F31").Value = 118
$ws.Range("F36").Value = 201
$ws.Range("F38").Value = 346
$ws.Range("F40").Value = 1058
$ws.Range("F42").Value = 77
$ws.Range("F43").Value = 2053
$ws.Range("F44").Value = 575
$ws.Range("F45").Value = 5
$ws.Range("F46").Value = 11
$ws.Range("F47").Value = 35
